# Crime_Data.xlsx edit: drop the 2011 literacy column from the "Literacy"
# sheet, keeping only the state name and the 2022 figure (old column C
# slides left into column B). This also removes the now-orphaned shared
# string "-" (it was only used as a placeholder in the deleted 2011
# column), which is why every shared-string index used elsewhere in the
# workbook shifts down by one automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Literacy")

# Delete column B ("2011"); column C ("2022") becomes the new column B.
[void]$ws.Columns.Item(2).Delete()

# Literacy becomes the active sheet / tab, with column B selected.
$ws.Activate() | Out-Null
$ws.Range("B1:B1048576").Select() | Out-Null
